$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 2 (the existing reference/template row) onto row 3
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)

# Match row 3's height to row 2's (18.75 -> 31.5)
$ws.Rows.Item(3).RowHeight = 31.5

# Fill in the new transaction data for row 3
$ws.Range("A3").Value = "18.09.2020"
$ws.Range("C3").Value = -100
$ws.Range("D3").Value = "UAH"
$ws.Range("E3").Value = 28
$ws.Range("F3").Formula = "=(B3+C3)/E3"
$ws.Range("G3").Value = "CARD"
$ws.Range("H3").Value = "Payment for transition action schemes from paper-view to draw.io"

# Update the dependent formulas so their cached values refresh
$ws.Range("I3").Formula = "=F3/2"
$ws.Range("J3").Formula = "=K3-I3"
$ws.Range("K3").Formula = "=K2+F3"

# Move the active selection to D3
$ws.Range("D3").Select() | Out-Null
